# "Changes of Rate Verification"
#
# Column C (rows 2-22) holds FedEx tracking numbers as TEXT (shared
# strings) even though they look numeric. Rows 5,6,7,13,14,15,16,17 also
# mirror the same tracking number in column D. This script replaces each
# of those tracking numbers with a new one, keeping the cells as text
# (no numeric conversion) and without touching any cell styles.
#
# Note: assigning a numeric-looking string straight to Range.Value (even
# with a leading apostrophe) makes Excel coerce it to a real number, or
# tags the cell with a new "quote prefix" style - neither of which match
# the source data (plain shared-string text, original style untouched).
# Routing the text through a TEXT() formula on a scratch cell and then
# using PasteSpecial(xlPasteValues) copies the *text* result into the
# destination cell, preserving both the string type and original style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new tracking number for column C (and D where noted below)
$newValues = [ordered]@{
  2  = "320018714339"
  3  = "320018714340"
  4  = "320018714372"
  5  = "320018714394"
  6  = "320018714431"
  7  = "320018714453"
  8  = "320018714486"
  9  = "320018714501"
  10 = "320018714534"
  11 = "320018714556"
  12 = "320018714590"
  13 = "320018714615"
  14 = "320018714648"
  15 = "320018714660"
  16 = "320018714692"
  17 = "320018714718"
  18 = "320018714751"
  19 = "320018714773"
  20 = "320018714800"
  21 = "320018714821"
  22 = "320018714854"
}

# rows where column D carries the same tracking number as column C
$dMirrorRows = @(5, 6, 7, 13, 14, 15, 16, 17)

# scratch cell well outside the used range (A1:AE24), restored/cleared at the end
$scratch = $ws.Range("ZZ1000")

foreach ($r in $newValues.Keys) {
    $val = $newValues[$r]

    $scratch.Formula = '=TEXT(' + $val + ',"0")'
    $scratch.Copy()
    $ws.Range("C$r").PasteSpecial(-4163)  # xlPasteValues

    if ($dMirrorRows -contains $r) {
        $ws.Range("D$r").PasteSpecial(-4163)  # xlPasteValues
    }
}

$scratch.ClearContents()
$scratch.Clear()

Write-Output "Updated tracking numbers for C2:C22 (and mirrored D cells)."
